$d = $word.ActiveDocument

# 1. Remove the _GoBack bookmark from its current location (end of the
#    "...Sortino ratio." paragraph).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Locate the "In this reference, ..." paragraph and the empty paragraph
# that immediately precedes it (the one right after the numbered
# reference entry).
$refParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "^In this reference") {
        $refParaIndex = $i
        break
    }
}

# 2. Re-create the _GoBack bookmark as a zero-length bookmark at the end of
#    that preceding empty paragraph (i.e. right before its paragraph mark).
$target = $d.Paragraphs($refParaIndex - 1)
$bmRange = $target.Range
[void]$bmRange.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 3. Delete the trailing, mostly-empty paragraphs that follow the
#    "In this reference, ..." paragraph, right up to (but excluding) the
#    final section break.
$lastParaIndex = $d.Paragraphs.Count
if ($refParaIndex -lt $lastParaIndex) {
    $delStart = $d.Paragraphs($refParaIndex + 1).Range.Start
    $delEnd = $d.Paragraphs($lastParaIndex).Range.End
    $delRange = $d.Range($delStart, $delEnd)
    $delRange.Delete()
}
